# BDD_geometrie.xlsx - "Ajout de fonctions de manipulation des bases de donnees"
#
# The French accented labels used for the fitting-type column ("A") are
# normalised to their unaccented ASCII equivalents, a couple of extra
# column widths are fixed so the labels aren't truncated, and the view is
# scrolled / the active cell moved down onto the data table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row --------------------------------------------------------
$ws.Range("D1").Value2 = "rapport diametre sortie entree"

# --- "te" (was "te with an accent") rows --------------------------------
for ($r = 7; $r -le 11; $r++) {
    $ws.Cells.Item($r, 1).Value2 = "te"
}

# --- "retrecissement brusque" rows --------------------------------------
for ($r = 28; $r -le 32; $r++) {
    $ws.Cells.Item($r, 1).Value2 = "retrecissement brusque"
}

# --- "retrecissement" rows -----------------------------------------------
for ($r = 33; $r -le 38; $r++) {
    $ws.Cells.Item($r, 1).Value2 = "retrecissement"
}

# --- Column widths so the longer labels are fully visible ---------------
$ws.Columns.Item(2).ColumnWidth = 25.8
$ws.Columns.Item(4).ColumnWidth = 24.95

# --- View / selection state ---------------------------------------------
# (The engine always pins the active cell to the top-left corner of
# whatever was most recently selected/activated, so the richest
# reproducible state is "the whole table A1:E38 selected"; a real Excel
# session additionally nudges the active cell down to A7 and scrolls the
# viewport to row 5, which we still issue for fidelity even though this
# host does not persist topLeftCell / a non-anchor active cell.)
$excel.ActiveWindow.ScrollRow = 5
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A7").Activate()
$ws.Range("A1:E38").Select()
